# database column types changed
#
# Slide 3 ("images") contains a table documenting the columns of a
# database table. The "iso" column's type changes from "double" to
# "int", and the "exposure_time" column's type changes from "string"
# to "int".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table

# Row 1 = header names, Row 2 = column data types.
# Column 8 = "iso" (was "double"), Column 10 = "exposure_time" (was "string")
$tbl.Cell(2, 8).Shape.TextFrame.TextRange.Text = "int"
$tbl.Cell(2, 10).Shape.TextFrame.TextRange.Text = "int"
